# 224613 Add --force-create and publish/review actions for items
#
# - General: new "Error" column header (C1), refreshed Product ID (PRD-),
#   selection moved to C3.
# - Parameters Groups / Items Groups / Agreements Parameters /
#   Item Parameters / Request Parameters / Subscription Parameters / Items /
#   Templates: refreshed generated IDs (9984-1895 -> 1213-3316).
# - Item Parameters: action list no longer offers "create", instead offers
#   "review"/"publish" (with trailing "--force-create" empty option), and the
#   selection moves to E13.
# - Items: new "Error" column header (S1), action list switches to
#   update/review/publish (no more create/delete), becomes the active sheet
#   with S2:S3 selected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# General
# ---------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("C1").Value = "Error"
$wsGeneral.Range("B3").Value = "PRD-1213-3316"
[void]$wsGeneral.Range("C3").Select()

# ---------------------------------------------------------------------
# Parameters Groups
# ---------------------------------------------------------------------
$wsParamGroups = $wb.Worksheets.Item("Parameters Groups")
$wsParamGroups.Range("A2").Value = "PGR-1213-3316-0002"
$wsParamGroups.Range("A3").Value = "PGR-1213-3316-0003"

# ---------------------------------------------------------------------
# Items Groups
# ---------------------------------------------------------------------
$wsItemGroups = $wb.Worksheets.Item("Items Groups")
$wsItemGroups.Range("A2").Value = "IGR-1213-3316-0002"
$wsItemGroups.Range("A3").Value = "IGR-1213-3316-0003"

# ---------------------------------------------------------------------
# Agreements Parameters
# ---------------------------------------------------------------------
$wsAgreementParams = $wb.Worksheets.Item("Agreements Parameters")
$wsAgreementParams.Range("A2").Value = "PAR-1213-3316-0001"
$wsAgreementParams.Range("H2").Value = "PGR-1213-3316-0002"
$wsAgreementParams.Range("A3").Value = "PAR-1213-3316-0002"

# ---------------------------------------------------------------------
# Item Parameters
# ---------------------------------------------------------------------
$wsItemParams = $wb.Worksheets.Item("Item Parameters")
$wsItemParams.Range("A2").Value = "PAR-1213-3316-0003"
$wsItemParams.Range("A3").Value = "PAR-1213-3316-0004"
$wsItemParams.Range("E2:E3").Validation.Formula1 = '"-,update,review,publish,"'
[void]$wsItemParams.Range("E13").Select()

# ---------------------------------------------------------------------
# Request Parameters
# ---------------------------------------------------------------------
$wsRequestParams = $wb.Worksheets.Item("Request Parameters")
$wsRequestParams.Range("A2").Value = "PAR-1213-3316-0005"
$wsRequestParams.Range("A3").Value = "PAR-1213-3316-0006"

# ---------------------------------------------------------------------
# Subscription Parameters
# ---------------------------------------------------------------------
$wsSubscriptionParams = $wb.Worksheets.Item("Subscription Parameters")
$wsSubscriptionParams.Range("A2").Value = "PAR-1213-3316-0007"
$wsSubscriptionParams.Range("A3").Value = "PAR-1213-3316-0008"

# ---------------------------------------------------------------------
# Templates
# ---------------------------------------------------------------------
$wsTemplates = $wb.Worksheets.Item("Templates")
$wsTemplates.Range("A2").Value = "TPL-1213-3316-0005"
$wsTemplates.Range("F2").Value = "Test content **Azure** {{ PAR-1213-3316-0001 }}"
$wsTemplates.Range("A3").Value = "TPL-1213-3316-0006"

# ---------------------------------------------------------------------
# Items (edited + activated last so it ends up the active tab)
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items")
$wsItems.Range("S1").Value = "Error"
$wsItems.Range("A2").Value = "ITM-1213-3316-0001"
$wsItems.Range("J2").Value = "IGR-1213-3316-0002"
$wsItems.Range("A3").Value = "ITM-1213-3316-0002"
$wsItems.Range("J3").Value = "IGR-1213-3316-0002"
$wsItems.Range("C2:C3").Validation.Formula1 = '"-,update,review,publish"'
[void]$wsItems.Range("S2:S3").Select()
